$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, "A").Value = "ECs"
$ws.Cells.Item(2, "B").Value = "Fn1"
$ws.Cells.Item(2, "C").Value = "Itgb6"
$ws.Cells.Item(2, "D").Value = "ECs"
$ws.Cells.Item(2, "E").Value = 3.0
$ws.Cells.Item(2, "F").Value = 1.0
$ws.Cells.Item(2, "G").Value = 40.75339133333333
$ws.Cells.Item(2, "H").Value = 122.260174
$ws.Cells.Item(2, "I").Value = 0.02126536631186857
$ws.Cells.Item(2, "J").Value = 0.02126536631186857
$ws.Cells.Item(2, "K").Value = 1.0
$ws.Cells.Item(2, "L").Value = 0.3333333333333333
$ws.Cells.Item(2, "M").Value = 0.111678
$ws.Cells.Item(2, "N").Value = 0.335034
$ws.Cells.Item(2, "O").Value = 0.01499007198665366
$ws.Cells.Item(2, "P").Value = 0.01499007198665366
$ws.Cells.Item(2, "Q").Value = 4.551257237324
$ws.Cells.Item(2, "R").Value = 40.961315135916
$ws.Cells.Item(2, "S").Value = 0.0003187693718374694
$ws.Cells.Item(2, "T").Value = 0.0003187693718374694

$ws.Cells.Item(3, "A").Value = "ECs"
$ws.Cells.Item(3, "B").Value = "Fn1"
$ws.Cells.Item(3, "C").Value = "Itgb6"
$ws.Cells.Item(3, "D").Value = "FAPs"
$ws.Cells.Item(3, "E").Value = 3.0
$ws.Cells.Item(3, "F").Value = 1.0
$ws.Cells.Item(3, "G").Value = 40.75339133333333
$ws.Cells.Item(3, "H").Value = 122.260174
$ws.Cells.Item(3, "I").Value = 0.02126536631186857
$ws.Cells.Item(3, "J").Value = 0.02126536631186857
$ws.Cells.Item(3, "K").Value = 3.0
$ws.Cells.Item(3, "L").Value = 1.0
$ws.Cells.Item(3, "M").Value = 2.242708666666667
$ws.Cells.Item(3, "N").Value = 6.728126
$ws.Cells.Item(3, "O").Value = 0.3010294270888212
$ws.Cells.Item(3, "P").Value = 0.3010294270888212
$ws.Cells.Item(3, "Q").Value = 91.39798393932489
$ws.Cells.Item(3, "R").Value = 822.581855453924
$ws.Cells.Item(3, "S").Value = 0.006401501037695714
$ws.Cells.Item(3, "T").Value = 0.006401501037695714

$ws.Cells.Item(4, "A").Value = "ECs"
$ws.Cells.Item(4, "B").Value = "Fn1"
$ws.Cells.Item(4, "C").Value = "Itgb6"
$ws.Cells.Item(4, "D").Value = "M2"
$ws.Cells.Item(4, "E").Value = 3.0
$ws.Cells.Item(4, "F").Value = 1.0
$ws.Cells.Item(4, "G").Value = 40.75339133333333
$ws.Cells.Item(4, "H").Value = 122.260174
$ws.Cells.Item(4, "I").Value = 0.02126536631186857
$ws.Cells.Item(4, "J").Value = 0.02126536631186857
$ws.Cells.Item(4, "K").Value = 1.0
$ws.Cells.Item(4, "L").Value = 0.3333333333333333
$ws.Cells.Item(4, "M").Value = 0.02011033333333333
$ws.Cells.Item(4, "N").Value = 0.060331
$ws.Cells.Item(4, "O").Value = 0.002699326137128775
$ws.Cells.Item(4, "P").Value = 0.002699326137128775
$ws.Cells.Item(4, "Q").Value = 0.8195642841771111
$ws.Cells.Item(4, "R").Value = 7.376078557594001
$ws.Cells.Item(4, "S").Value = 0.00005740215910124455
$ws.Cells.Item(4, "T").Value = 0.00005740215910124457

$ws.Cells.Item(5, "A").Value = "ECs"
$ws.Cells.Item(5, "B").Value = "Fn1"
$ws.Cells.Item(5, "C").Value = "Itgb6"
$ws.Cells.Item(5, "D").Value = "sCs"
$ws.Cells.Item(5, "E").Value = 3.0
$ws.Cells.Item(5, "F").Value = 1.0
$ws.Cells.Item(5, "G").Value = 40.75339133333333
$ws.Cells.Item(5, "H").Value = 122.260174
$ws.Cells.Item(5, "I").Value = 0.02126536631186857
$ws.Cells.Item(5, "J").Value = 0.02126536631186857
$ws.Cells.Item(5, "K").Value = 3.0
$ws.Cells.Item(5, "L").Value = 1.0
$ws.Cells.Item(5, "M").Value = 5.075634
$ws.Cells.Item(5, "N").Value = 15.226902
$ws.Cells.Item(5, "O").Value = 0.6812811747873964
$ws.Cells.Item(5, "P").Value = 0.6812811747873964
$ws.Cells.Item(5, "Q").Value = 206.849298666772
$ws.Cells.Item(5, "R").Value = 1861.643688000948
$ws.Cells.Item(5, "S").Value = 0.01448769374323414
$ws.Cells.Item(5, "T").Value = 0.01448769374323414

$ws.Cells.Item(6, "A").Value = "FAPs"
$ws.Cells.Item(6, "B").Value = "Fn1"
$ws.Cells.Item(6, "C").Value = "Itgb6"
$ws.Cells.Item(6, "D").Value = "ECs"
$ws.Cells.Item(6, "E").Value = 3.0
$ws.Cells.Item(6, "F").Value = 1.0
$ws.Cells.Item(6, "G").Value = 1689.289306666667
$ws.Cells.Item(6, "H").Value = 5067.86792
$ws.Cells.Item(6, "I").Value = 0.8814813868902838
$ws.Cells.Item(6, "J").Value = 0.8814813868902838
$ws.Cells.Item(6, "K").Value = 1.0
$ws.Cells.Item(6, "L").Value = 0.3333333333333333
$ws.Cells.Item(6, "M").Value = 0.111678
$ws.Cells.Item(6, "N").Value = 0.335034
$ws.Cells.Item(6, "O").Value = 0.01499007198665366
$ws.Cells.Item(6, "P").Value = 0.01499007198665366
$ws.Cells.Item(6, "Q").Value = 188.65645118992
$ws.Cells.Item(6, "R").Value = 1697.90806070928
$ws.Cells.Item(6, "S").Value = 0.01321346944438066
$ws.Cells.Item(6, "T").Value = 0.01321346944438066

$ws.Cells.Item(7, "A").Value = "FAPs"
$ws.Cells.Item(7, "B").Value = "Fn1"
$ws.Cells.Item(7, "C").Value = "Itgb6"
$ws.Cells.Item(7, "D").Value = "FAPs"
$ws.Cells.Item(7, "E").Value = 3.0
$ws.Cells.Item(7, "F").Value = 1.0
$ws.Cells.Item(7, "G").Value = 1689.289306666667
$ws.Cells.Item(7, "H").Value = 5067.86792
$ws.Cells.Item(7, "I").Value = 0.8814813868902838
$ws.Cells.Item(7, "J").Value = 0.8814813868902838
$ws.Cells.Item(7, "K").Value = 3.0
$ws.Cells.Item(7, "L").Value = 1.0
$ws.Cells.Item(7, "M").Value = 2.242708666666667
$ws.Cells.Item(7, "N").Value = 6.728126
$ws.Cells.Item(7, "O").Value = 0.3010294270888212
$ws.Cells.Item(7, "P").Value = 0.3010294270888212
$ws.Cells.Item(7, "Q").Value = 3788.583768568657
$ws.Cells.Item(7, "R").Value = 34097.25391711792
$ws.Cells.Item(7, "S").Value = 0.2653518368850417
$ws.Cells.Item(7, "T").Value = 0.2653518368850417

$ws.Cells.Item(8, "A").Value = "FAPs"
$ws.Cells.Item(8, "B").Value = "Fn1"
$ws.Cells.Item(8, "C").Value = "Itgb6"
$ws.Cells.Item(8, "D").Value = "M2"
$ws.Cells.Item(8, "E").Value = 3.0
$ws.Cells.Item(8, "F").Value = 1.0
$ws.Cells.Item(8, "G").Value = 1689.289306666667
$ws.Cells.Item(8, "H").Value = 5067.86792
$ws.Cells.Item(8, "I").Value = 0.8814813868902838
$ws.Cells.Item(8, "J").Value = 0.8814813868902838
$ws.Cells.Item(8, "K").Value = 1.0
$ws.Cells.Item(8, "L").Value = 0.3333333333333333
$ws.Cells.Item(8, "M").Value = 0.02011033333333333
$ws.Cells.Item(8, "N").Value = 0.060331
$ws.Cells.Item(8, "O").Value = 0.002699326137128775
$ws.Cells.Item(8, "P").Value = 0.002699326137128775
$ws.Cells.Item(8, "Q").Value = 33.97217105350222
$ws.Cells.Item(8, "R").Value = 305.74953948152
$ws.Cells.Item(8, "S").Value = 0.002379405747025465
$ws.Cells.Item(8, "T").Value = 0.002379405747025465

$ws.Cells.Item(9, "A").Value = "FAPs"
$ws.Cells.Item(9, "B").Value = "Fn1"
$ws.Cells.Item(9, "C").Value = "Itgb6"
$ws.Cells.Item(9, "D").Value = "sCs"
$ws.Cells.Item(9, "E").Value = 3.0
$ws.Cells.Item(9, "F").Value = 1.0
$ws.Cells.Item(9, "G").Value = 1689.289306666667
$ws.Cells.Item(9, "H").Value = 5067.86792
$ws.Cells.Item(9, "I").Value = 0.8814813868902838
$ws.Cells.Item(9, "J").Value = 0.8814813868902838
$ws.Cells.Item(9, "K").Value = 3.0
$ws.Cells.Item(9, "L").Value = 1.0
$ws.Cells.Item(9, "M").Value = 5.075634
$ws.Cells.Item(9, "N").Value = 15.226902
$ws.Cells.Item(9, "O").Value = 0.6812811747873964
$ws.Cells.Item(9, "P").Value = 0.6812811747873964
$ws.Cells.Item(9, "Q").Value = 8574.21424075376
$ws.Cells.Item(9, "R").Value = 77167.92816678384
$ws.Cells.Item(9, "S").Value = 0.600536674813836
$ws.Cells.Item(9, "T").Value = 0.600536674813836

$ws.Cells.Item(10, "A").Value = "M2"
$ws.Cells.Item(10, "B").Value = "Fn1"
$ws.Cells.Item(10, "C").Value = "Itgb6"
$ws.Cells.Item(10, "D").Value = "ECs"
$ws.Cells.Item(10, "E").Value = 3.0
$ws.Cells.Item(10, "F").Value = 1.0
$ws.Cells.Item(10, "G").Value = 100.9654023333333
$ws.Cells.Item(10, "H").Value = 302.896207
$ws.Cells.Item(10, "I").Value = 0.05268435816499466
$ws.Cells.Item(10, "J").Value = 0.05268435816499466
$ws.Cells.Item(10, "K").Value = 1.0
$ws.Cells.Item(10, "L").Value = 0.3333333333333333
$ws.Cells.Item(10, "M").Value = 0.111678
$ws.Cells.Item(10, "N").Value = 0.335034
$ws.Cells.Item(10, "O").Value = 0.01499007198665366
$ws.Cells.Item(10, "P").Value = 0.01499007198665366
$ws.Cells.Item(10, "Q").Value = 11.275614201782
$ws.Cells.Item(10, "R").Value = 101.480527816038
$ws.Cells.Item(10, "S").Value = 0.0007897423214639143
$ws.Cells.Item(10, "T").Value = 0.0007897423214639144

$ws.Cells.Item(11, "A").Value = "M2"
$ws.Cells.Item(11, "B").Value = "Fn1"
$ws.Cells.Item(11, "C").Value = "Itgb6"
$ws.Cells.Item(11, "D").Value = "FAPs"
$ws.Cells.Item(11, "E").Value = 3.0
$ws.Cells.Item(11, "F").Value = 1.0
$ws.Cells.Item(11, "G").Value = 100.9654023333333
$ws.Cells.Item(11, "H").Value = 302.896207
$ws.Cells.Item(11, "I").Value = 0.05268435816499466
$ws.Cells.Item(11, "J").Value = 0.05268435816499466
$ws.Cells.Item(11, "K").Value = 3.0
$ws.Cells.Item(11, "L").Value = 1.0
$ws.Cells.Item(11, "M").Value = 2.242708666666667
$ws.Cells.Item(11, "N").Value = 6.728126
$ws.Cells.Item(11, "O").Value = 0.3010294270888212
$ws.Cells.Item(11, "P").Value = 0.3010294270888212
$ws.Cells.Item(11, "Q").Value = 226.4359828464536
$ws.Cells.Item(11, "R").Value = 2037.923845618082
$ws.Cells.Item(11, "S").Value = 0.0158595421549506
$ws.Cells.Item(11, "T").Value = 0.0158595421549506

$ws.Cells.Item(12, "A").Value = "M2"
$ws.Cells.Item(12, "B").Value = "Fn1"
$ws.Cells.Item(12, "C").Value = "Itgb6"
$ws.Cells.Item(12, "D").Value = "M2"
$ws.Cells.Item(12, "E").Value = 3.0
$ws.Cells.Item(12, "F").Value = 1.0
$ws.Cells.Item(12, "G").Value = 100.9654023333333
$ws.Cells.Item(12, "H").Value = 302.896207
$ws.Cells.Item(12, "I").Value = 0.05268435816499466
$ws.Cells.Item(12, "J").Value = 0.05268435816499466
$ws.Cells.Item(12, "K").Value = 1.0
$ws.Cells.Item(12, "L").Value = 0.3333333333333333
$ws.Cells.Item(12, "M").Value = 0.02011033333333333
$ws.Cells.Item(12, "N").Value = 0.060331
$ws.Cells.Item(12, "O").Value = 0.002699326137128775
$ws.Cells.Item(12, "P").Value = 0.002699326137128775
$ws.Cells.Item(12, "Q").Value = 2.030447896057444
$ws.Cells.Item(12, "R").Value = 18.274031064517
$ws.Cells.Item(12, "S").Value = 0.0001422122650126238
$ws.Cells.Item(12, "T").Value = 0.0001422122650126239

$ws.Cells.Item(13, "A").Value = "M2"
$ws.Cells.Item(13, "B").Value = "Fn1"
$ws.Cells.Item(13, "C").Value = "Itgb6"
$ws.Cells.Item(13, "D").Value = "sCs"
$ws.Cells.Item(13, "E").Value = 3.0
$ws.Cells.Item(13, "F").Value = 1.0
$ws.Cells.Item(13, "G").Value = 100.9654023333333
$ws.Cells.Item(13, "H").Value = 302.896207
$ws.Cells.Item(13, "I").Value = 0.05268435816499466
$ws.Cells.Item(13, "J").Value = 0.05268435816499466
$ws.Cells.Item(13, "K").Value = 3.0
$ws.Cells.Item(13, "L").Value = 1.0
$ws.Cells.Item(13, "M").Value = 5.075634
$ws.Cells.Item(13, "N").Value = 15.226902
$ws.Cells.Item(13, "O").Value = 0.6812811747873964
$ws.Cells.Item(13, "P").Value = 0.6812811747873964
$ws.Cells.Item(13, "Q").Value = 512.463428906746
$ws.Cells.Item(13, "R").Value = 4612.170860160713
$ws.Cells.Item(13, "S").Value = 0.03589286142356752
$ws.Cells.Item(13, "T").Value = 0.03589286142356752

$ws.Cells.Item(14, "A").Value = "sCs"
$ws.Cells.Item(14, "B").Value = "Fn1"
$ws.Cells.Item(14, "C").Value = "Itgb6"
$ws.Cells.Item(14, "D").Value = "ECs"
$ws.Cells.Item(14, "E").Value = 3.0
$ws.Cells.Item(14, "F").Value = 1.0
$ws.Cells.Item(14, "G").Value = 85.41274733333334
$ws.Cells.Item(14, "H").Value = 256.238242
$ws.Cells.Item(14, "I").Value = 0.04456888863285297
$ws.Cells.Item(14, "J").Value = 0.04456888863285297
$ws.Cells.Item(14, "K").Value = 1.0
$ws.Cells.Item(14, "L").Value = 0.3333333333333333
$ws.Cells.Item(14, "M").Value = 0.111678
$ws.Cells.Item(14, "N").Value = 0.335034
$ws.Cells.Item(14, "O").Value = 0.01499007198665366
$ws.Cells.Item(14, "P").Value = 0.01499007198665366
$ws.Cells.Item(14, "Q").Value = 9.538724796692001
$ws.Cells.Item(14, "R").Value = 85.848523170228
$ws.Cells.Item(14, "S").Value = 0.0006680908489716159
$ws.Cells.Item(14, "T").Value = 0.000668090848971616

$ws.Cells.Item(15, "A").Value = "sCs"
$ws.Cells.Item(15, "B").Value = "Fn1"
$ws.Cells.Item(15, "C").Value = "Itgb6"
$ws.Cells.Item(15, "D").Value = "FAPs"
$ws.Cells.Item(15, "E").Value = 3.0
$ws.Cells.Item(15, "F").Value = 1.0
$ws.Cells.Item(15, "G").Value = 85.41274733333334
$ws.Cells.Item(15, "H").Value = 256.238242
$ws.Cells.Item(15, "I").Value = 0.04456888863285297
$ws.Cells.Item(15, "J").Value = 0.04456888863285297
$ws.Cells.Item(15, "K").Value = 3.0
$ws.Cells.Item(15, "L").Value = 1.0
$ws.Cells.Item(15, "M").Value = 2.242708666666667
$ws.Cells.Item(15, "N").Value = 6.728126
$ws.Cells.Item(15, "O").Value = 0.3010294270888212
$ws.Cells.Item(15, "P").Value = 0.3010294270888212
$ws.Cells.Item(15, "Q").Value = 191.5559086882769
$ws.Cells.Item(15, "R").Value = 1724.003178194492
$ws.Cells.Item(15, "S").Value = 0.01341654701113321
$ws.Cells.Item(15, "T").Value = 0.01341654701113321

$ws.Cells.Item(16, "A").Value = "sCs"
$ws.Cells.Item(16, "B").Value = "Fn1"
$ws.Cells.Item(16, "C").Value = "Itgb6"
$ws.Cells.Item(16, "D").Value = "M2"
$ws.Cells.Item(16, "E").Value = 3.0
$ws.Cells.Item(16, "F").Value = 1.0
$ws.Cells.Item(16, "G").Value = 85.41274733333334
$ws.Cells.Item(16, "H").Value = 256.238242
$ws.Cells.Item(16, "I").Value = 0.04456888863285297
$ws.Cells.Item(16, "J").Value = 0.04456888863285297
$ws.Cells.Item(16, "K").Value = 1.0
$ws.Cells.Item(16, "L").Value = 0.3333333333333333
$ws.Cells.Item(16, "M").Value = 0.02011033333333333
$ws.Cells.Item(16, "N").Value = 0.060331
$ws.Cells.Item(16, "O").Value = 0.002699326137128775
$ws.Cells.Item(16, "P").Value = 0.002699326137128775
$ws.Cells.Item(16, "Q").Value = 1.717678819789111
$ws.Cells.Item(16, "R").Value = 15.459109378102
$ws.Cells.Item(16, "S").Value = 0.0001203059659894415
$ws.Cells.Item(16, "T").Value = 0.0001203059659894416

$ws.Cells.Item(17, "A").Value = "sCs"
$ws.Cells.Item(17, "B").Value = "Fn1"
$ws.Cells.Item(17, "C").Value = "Itgb6"
$ws.Cells.Item(17, "D").Value = "sCs"
$ws.Cells.Item(17, "E").Value = 3.0
$ws.Cells.Item(17, "F").Value = 1.0
$ws.Cells.Item(17, "G").Value = 85.41274733333334
$ws.Cells.Item(17, "H").Value = 256.238242
$ws.Cells.Item(17, "I").Value = 0.04456888863285297
$ws.Cells.Item(17, "J").Value = 0.04456888863285297
$ws.Cells.Item(17, "K").Value = 3.0
$ws.Cells.Item(17, "L").Value = 1.0
$ws.Cells.Item(17, "M").Value = 5.075634
$ws.Cells.Item(17, "N").Value = 15.226902
$ws.Cells.Item(17, "O").Value = 0.6812811747873964
$ws.Cells.Item(17, "P").Value = 0.6812811747873964
$ws.Cells.Item(17, "Q").Value = 433.5238443984761
$ws.Cells.Item(17, "R").Value = 3901.714599586284
$ws.Cells.Item(17, "S").Value = 0.03036394480675871
$ws.Cells.Item(17, "T").Value = 0.03036394480675871

Write-Output "done"